$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.296878333333333
$ws.Range("H2").Value = 3.890635
$ws.Range("I2").Value = 0.01774073260139904
$ws.Range("J2").Value = 0.02506266560199287
$ws.Range("M2").Value = 1.565239
$ws.Range("N2").Value = 4.695717
$ws.Range("Q2").Value = 2.029924545588333
$ws.Range("R2").Value = 18.269320910295
$ws.Range("S2").Value = 0.01774073260139904
$ws.Range("T2").Value = 0.02506266560199287

# Row 3
$ws.Range("I3").Value = 0.09943605305674341
$ws.Range("J3").Value = 0.1404751766759988
$ws.Range("M3").Value = 1.565239
$ws.Range("N3").Value = 4.695717
$ws.Range("Q3").Value = 11.377640898572
$ws.Range("R3").Value = 102.398768087148
$ws.Range("S3").Value = 0.09943605305674341
$ws.Range("T3").Value = 0.1404751766759988

# Row 4
$ws.Range("G4").Value = 0.09795233333333332
$ws.Range("H4").Value = 0.293857
$ws.Range("I4").Value = 0.001339945397100812
$ws.Range("J4").Value = 0.001892965987764162
$ws.Range("M4").Value = 1.565239
$ws.Range("N4").Value = 4.695717
$ws.Range("Q4").Value = 0.1533188122743333
$ws.Range("R4").Value = 1.379869310469
$ws.Range("S4").Value = 0.001339945397100812
$ws.Range("T4").Value = 0.001892965987764162

# Row 5
$ws.Range("G5").Value = 64.06892400000001
$ws.Range("H5").Value = 128.137848
$ws.Range("I5").Value = 0.8764350668284411
$ws.Range("J5").Value = 0.8254375019458241
$ws.Range("M5").Value = 1.565239
$ws.Range("N5").Value = 4.695717
$ws.Range("Q5").Value = 100.283178532836
$ws.Range("R5").Value = 601.6990711970161
$ws.Range("S5").Value = 0.8764350668284411
$ws.Range("T5").Value = 0.8254375019458241

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.3690323333333334
$ws.Range("H6").Value = 1.107097
$ws.Range("I6").Value = 0.005048202116315478
$ws.Range("J6").Value = 0.007131689788420014
$ws.Range("M6").Value = 1.565239
$ws.Range("N6").Value = 4.695717
$ws.Range("Q6").Value = 0.5776238003943334
$ws.Range("R6").Value = 5.198614203549
$ws.Range("S6").Value = 0.005048202116315478
$ws.Range("T6").Value = 0.007131689788420014
